# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
# New data lives in columns AD, AE, AF (30, 31, 32) right after the existing
# "Unnamed: 28" column (AC / 29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column titles
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold / centered / bordered header style used by the other
# header cells (e.g. A1) by copying just the formatting over.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row (2-46) gets the same team season record.
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 67
    $ws.Cells.Item($row, 31).Value = 95
    $ws.Cells.Item($row, 32).Value = 0
}
